# ADD results from server
# Update inverter capacity result values on the "2025", "2030", and "2035"
# sheets (row 2 of each sheet holds the single data record).

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item("2025")
$ws.Range("A2").Value = 0
$ws.Range("E2").Value = 0.3806676774358236
$ws.Range("G2").Value = 0.2494892361374995
$ws.Range("I2").Value = 0.3505002166666667
$ws.Range("L2").Value = 0.6250375
$ws.Range("M2").Value = 0.07807858333333334
$ws.Range("N2").Value = 12.59325462450016
$ws.Range("O2").Value = 3.116032793548749

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.04468706386250129
$ws.Range("E2").Value = 0.2352331883523505
$ws.Range("I2").Value = 0.2500038147053857
$ws.Range("L2").Value = 0.2772243519612813
$ws.Range("M2").Value = 0.04511633333333334
$ws.Range("N2").Value = 5.345427060255322
$ws.Range("O2").Value = 2.482704542036897

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.06687192943954244
$ws.Range("B2").Value = 0.03189452449816901
$ws.Range("E2").Value = 0.1897417762144118
$ws.Range("I2").Value = 0.4857584928572489
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04542566666666664
$ws.Range("N2").Value = 8.365157032394315
$ws.Range("O2").Value = 5.022989240847572
